# "Added Online status" - append a new localization row (code ONLINE) to
# the Sheet1 language table, mirroring the existing Code/English/Vietnamese
# layout used by every other row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = $ws.Cells.Item(39, 1).Row

$ws.Cells.Item($newRow, 1).Value = "ONLINE"
$ws.Cells.Item($newRow, 2).Value = "PLAYERS ONLINE: "
$ws.Cells.Item($newRow, 3).Value = "ĐANG ONLINE: "

# Keep the pre-existing selection on C39 (last populated cell), matching
# the author's workbook state after the edit.
[void]$ws.Range("C39").Select()
